# "more intro stuff, implement intro"
# Insert 7 new localization rows (intro dialogue / attack-blob strings) right
# before the existing "lesson1_intro_1" row, shifting all subsequent rows
# down by 7 (old row 35 -> new row 42 ... old row 74 -> new row 81).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows at row 35 (pushes old rows 35-74 down to 42-81).
$ws.Range("A35:A41").EntireRow.Insert()

# Row insert copies formatting from the row above, which would create stray
# styles once we mix in wrapText / vertical-center further down - strip any
# inherited formatting first so the new block starts from the default style.
$ws.Range("A35:B41").ClearFormats()

# New Key / Value pairs. Written in this exact cell order so the shared
# string table (xl/sharedStrings.xml) grows with the same de-duplicated
# ordering the original author's edit produced.
$ws.Range("A35").Value = "intro_attack_blob"
$ws.Range("B35").Value = "Attack Blob"

$ws.Range("A36").Value = "intro_dialog_1"
$ws.Range("A37").Value = "intro_dialog_2"
$ws.Range("A38").Value = "intro_dialog_3"

$ws.Range("B36").Value = "Multiple space blobs have pierced through our dimension!"
$ws.Range("B37").Value = "Emergency protocol initiated."
$ws.Range("B38").Value = "We must banish them immediately before they fall down to Earth!"

$ws.Range("A39").Value = "intro_attack_1"
$ws.Range("A40").Value = "intro_attack_2"
$ws.Range("A41").Value = "intro_attack_3"

$ws.Range("B39").Value = "With our latest advancements in blobology, we will be deploying Attack Blobs."
$ws.Range("B41").Value = "Our intrepid hero, go forth, and use your mathematical might to banish these invading blobs back to their dimension!"
$ws.Range("B40").Value = "These blobs must be made with the power of multiplication, and who better to do it than you!"

# Matches the original sheet's formatting conventions: wrapped text on the
# short "Attack Blob" label, vertical-centered text on the row that ends up
# next to the taller wrapped rows.
$ws.Range("B35").WrapText = $true
$ws.Range("B40").VerticalAlignment = -4108

# Move the selection to reflect where editing left off.
$ws.Range("B40").Select()
